$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$values = @(5, 5, 0, 0, 12, 5, 7, 8, 2, 0)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.NumberFormat = "#,##0"
    $cell.Value = $values[$i]
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}
